# Apply cryptos list price/volume update (auto-generated from upstream diff)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.997.26"
$ws.Range("E2").Value = "  -1.23%  "
$ws.Range("D3").Value = "3.097.21"
$ws.Range("E3").Value = "  -1.39%  "
$ws.Range("E4").Value = "  -0.84%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "592.49"
$ws.Range("E5").Value = "  +1.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.48"
$ws.Range("E6").Value = "  +3.09%  "
$ws.Range("E7").Value = "  -2.16%  "
$ws.Range("E8").Value = "  +1.07%  "
$ws.Range("D9").Value = "3.096.77"
$ws.Range("E9").Value = "  -1.35%  "
$ws.Range("E10").Value = "  -1.94%  "
$ws.Range("E11").Value = "  -0.40%  "
$ws.Range("E12").Value = "  -2.53%  "
$ws.Range("E13").Value = "  -3.34%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.95"
$ws.Range("E14").Value = "  -3.97%  "
$ws.Range("E15").Value = "  -0.79%  "
$ws.Range("D16").Value = "3.608.75"
$ws.Range("E16").Value = "  -1.72%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.20"
$ws.Range("E17").Value = "  -1.01%  "
$ws.Range("D18").Value = "63.910.29"
$ws.Range("E18").Value = "  -0.81%  "
$ws.Range("D19").Value = "3.100.47"
$ws.Range("E19").Value = "  -13.86%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "480.96"
$ws.Range("E20").Value = "  +1.65%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.46"
$ws.Range("E21").Value = "  -2.85%  "
$ws.Range("E22").Value = "  -5.59%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.56"
$ws.Range("E23").Value = "  -1.64%  "
$ws.Range("E24").Value = "  +3.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.41"
$ws.Range("E25").Value = "  -1.46%  "
$ws.Range("E26").Value = "  -3.99%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.73"
$ws.Range("E27").Value = "  +8.58%  "
$ws.Range("E28").Value = "  +0.17%  "
$ws.Range("E29").Value = "  +3.38%  "
$ws.Range("E30").Value = "  -1.40%  "
$ws.Range("E31").Value = "  -0.91%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.19"
$ws.Range("E32").Value = "  -1.25%  "
$ws.Range("E33").Value = "  -3.74%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "27.24"
$ws.Range("E34").Value = "  -1.92%  "
$ws.Range("D35").Value = "0.0₃0842"
$ws.Range("E35").Value = "  -3.66%  "
$ws.Range("E36").Value = "  -0.47%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.03"
$ws.Range("E37").Value = "  -2.79%  "
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.26"
$ws.Range("E38").Value = "  -2.21%  "
$ws.Range("B39").Value = "dogwifhat"
$ws.Range("C39").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.29"
$ws.Range("E39").Value = "  -5.66%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "51.04"
$ws.Range("E40").Value = "  -0.36%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.27"
$ws.Range("E41").Value = "  -0.54%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "442.63"
$ws.Range("E42").Value = "  -4.30%  "
$ws.Range("E43").Value = "  -2.21%  "
$ws.Range("E44").Value = "  +2.12%  "
$ws.Range("E45").Value = "  -2.61%  "
$ws.Range("E46").Value = "  +5.22%  "
$ws.Range("D47").Value = "2.835.50"
$ws.Range("E47").Value = "  -1.90%  "
$ws.Range("E48").Value = "  -0.10%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "26.11"
$ws.Range("E49").Value = "  +2.23%  "
$ws.Range("E50").Value = "  +0.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.24"
$ws.Range("E51").Value = "  -2.36%  "
